# Update the "想去人数" (want-to-go count) figures that changed between
# scrapes, on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 34
$ws1.Range("F5").Value = 365
$ws1.Range("F6").Value = 1938
$ws1.Range("F7").Value = 104

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 34
$ws4.Range("F5").Value = 365
$ws4.Range("F10").Value = 1938
$ws4.Range("F11").Value = 104
